# IKD update: GaN CMOS 2026-02-04T23:32Z
# Appends 8 new literature rows (128-135) to the Master sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column S holds ISO-formatted date strings ("AddedDate") that must stay literal
# text (matching the rest of the sheet) instead of being auto-parsed into Excel
# date serials, so force Text format on that column for the new rows up front.
$ws.Range("S128:S135").NumberFormat = "@"

# Row 128
$ws.Cells.Item(128, 2).Value = 'Synthesis of Sliding Mode Control Strategy for T-Type Grid Inverter in Presence Grid Voltage Disturbance'
$ws.Cells.Item(128, 3).Value = 2026
$ws.Cells.Item(128, 4).Value = 'MDPI AG'
$ws.Cells.Item(128, 5).Value = 'Energies'
$ws.Cells.Item(128, 6).Value = 'Sawiński, Albert; Chudzik, Piotr; Tatar, Karol'
$ws.Cells.Item(128, 8).Value = '10.3390/en19030790'
$ws.Cells.Item(128, 9).Value = 'https://doi.org/10.3390/en19030790'
$ws.Cells.Item(128, 10).Value = 'Journal'
$ws.Cells.Item(128, 11).Value = 'Inverter'
$ws.Cells.Item(128, 12).Value = 'Experiment'
$ws.Cells.Item(128, 13).Value = 'Contacts'
$ws.Cells.Item(128, 17).Value = 'Synthesis of Sliding Mode Control Strategy for T-Type Grid Inverter in Presence Grid Voltage Disturbance'
$ws.Cells.Item(128, 18).Value = 'High'
$ws.Cells.Item(128, 19).Value = '2026-02-04'

# Row 129
$ws.Cells.Item(129, 2).Value = 'A comparative review of impact-induced failure mechanisms in monolithic and hybrid structures'
$ws.Cells.Item(129, 3).Value = 2026
$ws.Cells.Item(129, 4).Value = 'Elsevier BV'
$ws.Cells.Item(129, 5).Value = 'Results in Engineering'
$ws.Cells.Item(129, 6).Value = 'Zhong, Zhenhang; Wang, Jiacong'
$ws.Cells.Item(129, 8).Value = '10.1016/j.rineng.2026.109340'
$ws.Cells.Item(129, 9).Value = 'https://doi.org/10.1016/j.rineng.2026.109340'
$ws.Cells.Item(129, 10).Value = 'Journal'
$ws.Cells.Item(129, 11).Value = 'Co-integration'
$ws.Cells.Item(129, 12).Value = 'Hybrid'
$ws.Cells.Item(129, 13).Value = 'Integration'
$ws.Cells.Item(129, 17).Value = 'A comparative review of impact-induced failure mechanisms in monolithic and hybrid structures'
$ws.Cells.Item(129, 18).Value = 'High'
$ws.Cells.Item(129, 19).Value = '2026-02-04'

# Row 130
$ws.Cells.Item(130, 2).Value = 'Highly Sensitive Room-Temperature Graphene-Modulated AlGaN/GaN HEMT THz Detector Architecture'
$ws.Cells.Item(130, 3).Value = 2026
$ws.Cells.Item(130, 4).Value = 'MDPI AG'
$ws.Cells.Item(130, 5).Value = 'Sensors'
$ws.Cells.Item(130, 6).Value = 'Sengupta, Rudrarup; Sarusi, Gabby'
$ws.Cells.Item(130, 8).Value = '10.3390/s26031006'
$ws.Cells.Item(130, 9).Value = 'https://doi.org/10.3390/s26031006'
$ws.Cells.Item(130, 10).Value = 'Journal'
$ws.Cells.Item(130, 11).Value = 'n-FET'
$ws.Cells.Item(130, 12).Value = 'Experiment'
$ws.Cells.Item(130, 13).Value = 'Contacts'
$ws.Cells.Item(130, 17).Value = 'Highly Sensitive Room-Temperature Graphene-Modulated AlGaN/GaN HEMT THz Detector Architecture'
$ws.Cells.Item(130, 18).Value = 'High'
$ws.Cells.Item(130, 19).Value = '2026-02-04'

# Row 131
$ws.Cells.Item(131, 2).Value = 'Highly Sensitive Room-Temperature Graphene-Modulated AlGaN/GaN HEMT THz Detector Architecture'
$ws.Cells.Item(131, 3).Value = 2026
$ws.Cells.Item(131, 4).Value = 'MDPI AG'
$ws.Cells.Item(131, 5).Value = 'Sensors'
$ws.Cells.Item(131, 6).Value = 'Sengupta, Rudrarup; Sarusi, Gabby'
$ws.Cells.Item(131, 8).Value = '10.3390/s26031006'
$ws.Cells.Item(131, 9).Value = 'https://doi.org/10.3390/s26031006'
$ws.Cells.Item(131, 10).Value = 'Journal'
$ws.Cells.Item(131, 11).Value = 'n-FET'
$ws.Cells.Item(131, 12).Value = 'Experiment'
$ws.Cells.Item(131, 13).Value = 'Contacts'
$ws.Cells.Item(131, 17).Value = 'Highly Sensitive Room-Temperature Graphene-Modulated AlGaN/GaN HEMT THz Detector Architecture'
$ws.Cells.Item(131, 18).Value = 'High'
$ws.Cells.Item(131, 19).Value = '2026-02-04'

# Row 132
$ws.Cells.Item(132, 2).Value = 'Highly Sensitive Room-Temperature Graphene-Modulated AlGaN/GaN HEMT THz Detector Architecture'
$ws.Cells.Item(132, 3).Value = 2026
$ws.Cells.Item(132, 4).Value = 'MDPI AG'
$ws.Cells.Item(132, 5).Value = 'Sensors'
$ws.Cells.Item(132, 6).Value = 'Sengupta, Rudrarup; Sarusi, Gabby'
$ws.Cells.Item(132, 8).Value = '10.3390/s26031006'
$ws.Cells.Item(132, 9).Value = 'https://doi.org/10.3390/s26031006'
$ws.Cells.Item(132, 10).Value = 'Journal'
$ws.Cells.Item(132, 11).Value = 'n-FET'
$ws.Cells.Item(132, 12).Value = 'Experiment'
$ws.Cells.Item(132, 13).Value = 'Contacts'
$ws.Cells.Item(132, 17).Value = 'Highly Sensitive Room-Temperature Graphene-Modulated AlGaN/GaN HEMT THz Detector Architecture'
$ws.Cells.Item(132, 18).Value = 'High'
$ws.Cells.Item(132, 19).Value = '2026-02-04'

# Row 133
$ws.Cells.Item(133, 2).Value = 'Synthesis of Sliding Mode Control Strategy for T-Type Grid Inverter in Presence Grid Voltage Disturbance'
$ws.Cells.Item(133, 3).Value = 2026
$ws.Cells.Item(133, 4).Value = 'MDPI AG'
$ws.Cells.Item(133, 5).Value = 'Energies'
$ws.Cells.Item(133, 6).Value = 'Sawiński, Albert; Chudzik, Piotr; Tatar, Karol'
$ws.Cells.Item(133, 8).Value = '10.3390/en19030790'
$ws.Cells.Item(133, 9).Value = 'https://doi.org/10.3390/en19030790'
$ws.Cells.Item(133, 10).Value = 'Journal'
$ws.Cells.Item(133, 11).Value = 'Inverter'
$ws.Cells.Item(133, 12).Value = 'Experiment'
$ws.Cells.Item(133, 13).Value = 'Contacts'
$ws.Cells.Item(133, 17).Value = 'Synthesis of Sliding Mode Control Strategy for T-Type Grid Inverter in Presence Grid Voltage Disturbance'
$ws.Cells.Item(133, 18).Value = 'High'
$ws.Cells.Item(133, 19).Value = '2026-02-04'

# Row 134
$ws.Cells.Item(134, 2).Value = 'Assessment of an FPGA Implementation of a Hybrid PUF Based on a Configurable Transient Effect Ring Oscillator and Ring Oscillator (TERORO-PUF)'
$ws.Cells.Item(134, 3).Value = 2026
$ws.Cells.Item(134, 4).Value = 'MDPI AG'
$ws.Cells.Item(134, 5).Value = 'Electronics'
$ws.Cells.Item(134, 6).Value = 'Casado-Galán, Alejandro; Núñez, Juan; Tena-Sánchez, Erica; Potestad-Ordóñez, Francisco Eugenio; Acosta, Antonio José'
$ws.Cells.Item(134, 8).Value = '10.3390/electronics15030661'
$ws.Cells.Item(134, 9).Value = 'https://doi.org/10.3390/electronics15030661'
$ws.Cells.Item(134, 10).Value = 'Journal'
$ws.Cells.Item(134, 11).Value = 'Inverter'
$ws.Cells.Item(134, 12).Value = 'Experiment'
$ws.Cells.Item(134, 13).Value = 'Contacts'
$ws.Cells.Item(134, 17).Value = 'Assessment of an FPGA Implementation of a Hybrid PUF Based on a Configurable Transient Effect Ring Oscillator and Ring Oscillator (TERORO-PUF)'
$ws.Cells.Item(134, 18).Value = 'High'
$ws.Cells.Item(134, 19).Value = '2026-02-04'

# Row 135
$ws.Cells.Item(135, 2).Value = 'Synthesis of Sliding Mode Control Strategy for T-Type Grid Inverter in Presence Grid Voltage Disturbance'
$ws.Cells.Item(135, 3).Value = 2026
$ws.Cells.Item(135, 4).Value = 'MDPI AG'
$ws.Cells.Item(135, 5).Value = 'Energies'
$ws.Cells.Item(135, 6).Value = 'Sawiński, Albert; Chudzik, Piotr; Tatar, Karol'
$ws.Cells.Item(135, 8).Value = '10.3390/en19030790'
$ws.Cells.Item(135, 9).Value = 'https://doi.org/10.3390/en19030790'
$ws.Cells.Item(135, 10).Value = 'Journal'
$ws.Cells.Item(135, 11).Value = 'Inverter'
$ws.Cells.Item(135, 12).Value = 'Experiment'
$ws.Cells.Item(135, 13).Value = 'Contacts'
$ws.Cells.Item(135, 17).Value = 'Synthesis of Sliding Mode Control Strategy for T-Type Grid Inverter in Presence Grid Voltage Disturbance'
$ws.Cells.Item(135, 18).Value = 'High'
$ws.Cells.Item(135, 19).Value = '2026-02-04'
